$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 79 (shifts existing rows 79+ down by 3)
$ws.Rows("79:81").Insert()

# New row 79: NSCLC / v2.0-public / data_fusions / syn30381335 / 2022-10
$ws.Range("A79").Value = "NSCLC"
$ws.Range("B79").Value = "v2.0-public"
$ws.Range("C79").Value = "data_fusions"
$ws.Range("D79").Value = "syn30381335"
$ws.Range("E79").Value = "2022-10"

# New row 80: NSCLC / v2.0-public / data_cna / syn30381332 / 2022-10
$ws.Range("A80").Value = "NSCLC"
$ws.Range("B80").Value = "v2.0-public"
$ws.Range("C80").Value = "data_cna"
$ws.Range("D80").Value = "syn30381332"
$ws.Range("E80").Value = "2022-10"

# New row 81: NSCLC / v2.0-public / data_mutations_extended / syn30381327 / 2022-10
$ws.Range("A81").Value = "NSCLC"
$ws.Range("B81").Value = "v2.0-public"
$ws.Range("C81").Value = "data_mutations_extended"
$ws.Range("D81").Value = "syn30381327"
$ws.Range("E81").Value = "2022-10"

# Match the author's final selection / scroll position
$ws.Range("D81").Select()
$excel.ActiveWindow.ScrollRow = 67
